$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Extend header row 1 with two new columns (P=14, Q=15) ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Copy formatting (border/bold/alignment) from O1, the last existing header cell,
# onto the two freshly added header cells so they match the rest of the header row.
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# --- 2. Recalculated data rows 2-25 ---
# For every data row: columns C-I get new computed values, column O (the old
# "total" column) becomes 0, a new column P is added (always 0), and a new
# column Q is added holding the (updated) total that used to live in column O.
$rows = @{
  2 = @{ C=0.2062132122301676; D=0.03452013301317436; E=1.495232667404423; F=2.465111766229057; G=0.0007587804190388319; H=0.02803522330153196; I=0.07122987443985451; Q=8.625723999843501 }
  3 = @{ C=0.1825782504875662; D=0.03073011378493007; E=1.301123653223385; F=2.158698407357988; G=0.0007644417140412196; H=0.02137244059301224; I=0.05693383576570188; Q=7.559733382084062 }
  4 = @{ C=0.1678781834747554; D=0.02836975144844445; E=1.182777697230961; F=1.971577367073891; G=0.000768021576340372; H=0.01762511221425189; I=0.04874716355061981; Q=6.908178629627344 }
  5 = @{ C=0.1611424204649978; D=0.027322513805494; E=1.134645169277292; F=1.893663126126341; G=0.0007695229078780106; H=0.01616769533383411; I=0.04553924124969333; Q=6.636050624870506 }
  6 = @{ C=0.1591737534329098; D=0.02705438348821687; E=1.126572835850496; F=1.87846196071628; G=0.0007697929358507325; H=0.01592042019512108; I=0.04501749222994089; Q=6.582127798215652 }
  7 = @{ C=0.165458509939171; D=0.02809832874866913; E=1.181879590412478; F=1.964284057775515; G=0.000768093502567734; H=0.01757715677234839; I=0.04869291834943024; Q=6.880450964701936 }
  8 = @{ C=0.194965336559946; D=0.0328743248267287; E=1.427746137677929; F=2.350814355404452; G=0.0007607784515223119; H=0.02561700676978385; I=0.0661030887141445; Q=8.225074577546081 }
  9 = @{ C=0.2557785134384574; D=0.04246010601113426; E=1.919160950455648; F=3.131183301424244; G=0.0007471296242900018; H=0.04459914136076826; I=0.1058863026861863; Q=10.93917120733789 }
  10 = @{ C=0.2980539440814312; D=0.04793985769025255; E=2.180003320580923; F=3.679110950414156; G=0.0007378740379154656; H=0.05969667129688094; I=0.138558052099004; Q=12.82826607980957 }
  11 = @{ C=0.3076151637307305; D=0.03941282183964745; E=1.458506077811492; F=3.632020587140062; G=0.000736487558335665; H=0.07346859366028369; I=0.1447385136987958; Q=12.55067582091397 }
  12 = @{ C=0.3112181964743286; D=0.03237318057765393; E=0.9321154897466499; F=3.486940279146296; G=0.000736951370893036; H=0.1071589415353174; I=0.1434903718917848; Q=11.96096212953569 }
  13 = @{ C=0.3066640607967344; D=0.02570599905104487; E=0.5263433127709547; F=3.249504323903579; G=0.0007389244745406378; H=0.1573664906652539; I=0.1364432867728249; Q=11.0644344930916 }
  14 = @{ C=0.2991676459865857; D=0.02131805361349493; E=0.3164528817313297; F=3.042074817997502; G=0.0007409679880122742; H=0.2028147415036585; I=0.1291867130974191; Q=10.30041316872598 }
  15 = @{ C=0.2945543788243015; D=0.02017281835578544; E=0.2739382330171658; F=2.970077445361255; G=0.0007418333977551317; H=0.2139811819924518; I=0.1261702881499103; Q=10.04170610277822 }
  16 = @{ C=0.2761279780223589; D=0.01935092930526139; E=0.2635724612734762; F=2.786650979598932; G=0.0007453331103943323; H=0.196464873784322; I=0.1139313015565238; Q=9.427754723502915 }
  17 = @{ C=0.2658149304771626; D=0.02090291503434472; E=0.3522839429266824; F=2.754578388604102; G=0.0007469685434605837; H=0.1569452484738036; I=0.108314138216846; Q=9.356076777291037 }
  18 = @{ C=0.2635573844842867; D=0.02526202001759259; E=0.6015939595596933; F=2.854513301702525; G=0.00074702237257462; H=0.105485832718621; I=0.1077179054796789; Q=9.765487463008355 }
  19 = @{ C=0.2637408024227881; D=0.03195389919016023; E=1.063419647291042; F=3.045632383066703; G=0.0007457094319803367; H=0.06523954185914249; I=0.1117226668372275; Q=10.50496887608301 }
  20 = @{ C=0.2793407580425367; D=0.04562799013000074; E=2.105479798230888; F=3.512468536267107; G=0.0007404437960557399; H=0.05532613631724548; I=0.1293625464831027; Q=12.24524227118297 }
  21 = @{ C=0.3142969304566066; D=0.05196693224527849; E=2.477224327028608; F=3.993254688132936; G=0.0007329192563010345; H=0.07008815346526021; I=0.158039388381062; Q=13.92351275219653 }
  22 = @{ C=0.3400297655444717; D=0.05579569802971207; E=2.668247709860722; F=4.303050692851912; G=0.0007281613522351904; H=0.07980333804283823; I=0.1776531822115759; Q=15.00137203242127 }
  23 = @{ C=0.3292055944723273; D=0.05407481942569348; E=2.566455741422971; F=4.145268541673744; G=0.0007306388170778441; H=0.07461524531951014; I=0.1671655414525537; Q=14.45538754922916 }
  24 = @{ C=0.2837080296176708; D=0.04702718752439949; E=2.18577175694486; F=3.54277419665965; G=0.0007402221874129224; H=0.05627264951743793; I=0.1298794379051182; Q=12.36438140078747 }
  25 = @{ C=0.2350880966656774; D=0.03941135941203555; E=1.784415503561462; F=2.906734149520304; G=0.0007508329451554903; H=0.03898882851764895; I=0.09429756404644607; Q=10.15473619962268 }
}

foreach ($r in $rows.Keys) {
    $v = $rows[$r]
    $ws.Range("C$r").Value = $v.C
    $ws.Range("D$r").Value = $v.D
    $ws.Range("E$r").Value = $v.E
    $ws.Range("F$r").Value = $v.F
    $ws.Range("G$r").Value = $v.G
    $ws.Range("H$r").Value = $v.H
    $ws.Range("I$r").Value = $v.I
    $ws.Range("O$r").Value = 0
    $ws.Range("P$r").Value = 0
    $ws.Range("Q$r").Value = $v.Q
}
